$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new calibration parameters are appended to the table: subjective
# belief-shock standard deviations estimated from the Survey of Consumer
# Expectations (SCE), belonging to a new "subjective" block.

$ws.Range("A21").Value = "subjective"
$ws.Range("B21").Value = '$\sigma_\psi^{\text{sub}}$'
$ws.Range("C21").Value = 0.03457920401687286
$ws.Range("D21").Value = "estimated from SCE"

$ws.Range("A22").Value = "subjective"
$ws.Range("B22").Value = '$\sigma_\theta^{\text{sub}}$'
$ws.Range("C22").Value = 0.02010668171428303
$ws.Range("D22").Value = "estimated from SCE"

# Match the formatting convention used by the rest of the "block" column
# (bold, centered, bordered cell), same as the other cells in column A.
$ws.Range("A20").Copy()
$ws.Range("A21:A22").PasteSpecial(-4122)
$excel.CutCopyMode = 0
